$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = "ane981"
$ws.Range("E3").Value = "johnTheJOe121"
$ws.Range("E4").Value = "auraLAura901"
$ws.Range("E5").Value = "saratheS11"
$ws.Range("E6").Value = "lindaBone1231"
$ws.Range("E7").Value = "johnTrucker731"
